$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New VIN rows (4-9) added below the existing sample rows (2-3),
# part of the "symbol format" refactor (SYMBOL_2017 / BI00x / PD00x / UM00x / MP00x).

# 1) Clone the existing data-row formatting (style index 2, i.e. left-aligned)
#    from row 2 down onto the new rows so the new cells carry the same s="2" xf.
$ws.Range("A2:AJ2").Copy()
foreach ($r in 4..9) {
    $ws.Range("A" + $r + ":AJ" + $r).PasteSpecial(-4122)
}

# 2) Write the actual row values.
$newRows = @{
    4 = @{ "A"="1J2WW12P&5"; "B"="SYMBOL_2017"; "C"=2015; "D"="ACURA"; "E"="MDX"; "F"="Gt"; "G"="MDX ADVANCE"; "H"=53080; "I"="WAG"; "J"="UT_SS"; "K"="SUV"; "L"="UT_SS"; "M"="WAG"; "N"="4.5L V10"; "O"=8; "P"="G"; "Q"=214; "R"="2WD"; "S"=2; "T"="000R"; "U"="DUAL AIR BAGS FRONT"; "V"=2; "W"="4 WHEEL STANDARD"; "X"="STD"; "Y"="B-IMMOBILIZER/KEYLSS ENTRY/ALARM"; "Z"=42; "AA"=42; "AB"="Y"; "AC"="BI001"; "AD"="PD001"; "AE"="UM001"; "AF"="MP001"; "AG"=20010101; "AH"="Y"; "AI"="Y"; "AJ"="N" }
    5 = @{ "A"="1J2WW12P&5"; "B"="SYMBOL_2017"; "C"=2015; "D"="ACURA"; "E"="MDX"; "F"="Gt"; "G"="MDX ADVANCE"; "H"=53080; "I"="WAG"; "J"="UT_SS"; "K"="SUV"; "L"="UT_SS"; "M"="WAG"; "N"="4.5L V10"; "O"=8; "P"="G"; "Q"=214; "R"="2WD"; "S"=2; "T"="000R"; "U"="DUAL AIR BAGS FRONT"; "V"=2; "W"="4 WHEEL STANDARD"; "X"="STD"; "Y"="B-IMMOBILIZER/KEYLSS ENTRY/ALARM"; "Z"=42; "AA"=42; "AB"="Y"; "AC"="BI002"; "AD"="PD002"; "AE"="UM002"; "AF"="MP002"; "AG"=20000101; "AH"="N"; "AI"="Y"; "AJ"="N" }
    6 = @{ "A"="DDDKN3DD&E"; "B"="SYMBOL_2000"; "C"=2018; "D"="TOYOTA"; "E"="TOYOTA"; "F"="Gt"; "G"="MDX ADVANCE"; "H"=53080; "I"="WAG"; "J"="UT_SS"; "K"="SUV"; "L"="UT_SS"; "M"="WAG"; "N"="4.5L V10"; "O"=8; "P"="G"; "Q"=214; "R"="2WD"; "S"=2; "T"="000R"; "U"="DUAL AIR BAGS FRONT"; "V"=2; "W"="4 WHEEL STANDARD"; "X"="STD"; "Y"="B-IMMOBILIZER/KEYLSS ENTRY/ALARM"; "Z"=42; "AA"=42; "AB"="Y"; "AC"="C"; "AD"="C"; "AE"="C"; "AF"="C"; "AG"=20010101; "AH"="Y"; "AI"="Y"; "AJ"="N" }
    7 = @{ "A"="DDDKN3DD&E"; "B"="SYMBOL_2000"; "C"=2018; "D"="TOYOTA"; "E"="TOYOTA"; "F"="Gt"; "G"="MDX ADVANCE"; "H"=53080; "I"="WAG"; "J"="UT_SS"; "K"="SUV"; "L"="UT_SS"; "M"="WAG"; "N"="4.5L V10"; "O"=8; "P"="G"; "Q"=214; "R"="2WD"; "S"=2; "T"="000R"; "U"="DUAL AIR BAGS FRONT"; "V"=2; "W"="4 WHEEL STANDARD"; "X"="STD"; "Y"="B-IMMOBILIZER/KEYLSS ENTRY/ALARM"; "Z"=42; "AA"=42; "AB"="Y"; "AC"="N"; "AD"="N"; "AE"="N"; "AF"="N"; "AG"=20000101; "AH"="N"; "AI"="Y"; "AJ"="N" }
    8 = @{ "A"="DDDKN3DD&E"; "B"="SYMBOL_2017"; "C"=2018; "D"="TOYOTA"; "E"="TOYOTA"; "F"="Gt"; "G"="MDX ADVANCE"; "H"=53080; "I"="WAG"; "J"="UT_SS"; "K"="SUV"; "L"="UT_SS"; "M"="WAG"; "N"="4.5L V10"; "O"=8; "P"="G"; "Q"=214; "R"="2WD"; "S"=2; "T"="000R"; "U"="DUAL AIR BAGS FRONT"; "V"=2; "W"="4 WHEEL STANDARD"; "X"="STD"; "Y"="B-IMMOBILIZER/KEYLSS ENTRY/ALARM"; "Z"=42; "AA"=42; "AB"="Y"; "AC"="BI001"; "AD"="PD001"; "AE"="UM001"; "AF"="MP001"; "AG"=20010101; "AH"="Y"; "AI"="Y"; "AJ"="N" }
    9 = @{ "A"="DDDKN3DD&E"; "B"="SYMBOL_2017"; "C"=2018; "D"="TOYOTA"; "E"="TOYOTA"; "F"="Gt"; "G"="MDX ADVANCE"; "H"=53080; "I"="WAG"; "J"="UT_SS"; "K"="SUV"; "L"="UT_SS"; "M"="WAG"; "N"="4.5L V10"; "O"=8; "P"="G"; "Q"=214; "R"="2WD"; "S"=2; "T"="000R"; "U"="DUAL AIR BAGS FRONT"; "V"=2; "W"="4 WHEEL STANDARD"; "X"="STD"; "Y"="B-IMMOBILIZER/KEYLSS ENTRY/ALARM"; "Z"=42; "AA"=42; "AB"="Y"; "AC"="BI002"; "AD"="PD002"; "AE"="UM002"; "AF"="MP002"; "AG"=20000101; "AH"="N"; "AI"="Y"; "AJ"="N" }
}

foreach ($r in $newRows.Keys) {
    $rowVals = $newRows[$r]
    foreach ($col in $rowVals.Keys) {
        $ws.Range($col + $r).Value = $rowVals[$col]
    }
}

# 3) Restore the cursor position left by the edit (matches the saved sheetView selection).
$ws.Range("E17").Select()
